# chore: update Sheets via scheduled runner
# Refresh scraped market-price-derived columns (currentAveragePrice / NQ / HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across the per-job "Ragnarok_Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 92275.38
$ws.Range("J74").Value = 10000
$ws.Range("L74").Value = 10000
$ws.Range("N74").Value = -11872

$ws.Range("H77").Value = 92275.38
$ws.Range("J77").Value = 10000
$ws.Range("L77").Value = 50000
$ws.Range("N77").Value = -59360

$ws.Range("H92").Value = 2884.2273
$ws.Range("I92").Value = 1569.625
$ws.Range("J92").Value = 6389.8335
$ws.Range("K92").Value = 1569.625
$ws.Range("L92").Value = 6389.8335
$ws.Range("M92").Value = -321.625
$ws.Range("N92").Value = -8885.833500000001

$ws.Range("H94").Value = 3057.2307
$ws.Range("I94").Value = 3057.2307
$ws.Range("K94").Value = 3057.2307
$ws.Range("M94").Value = -2606.2307

$ws.Range("H100").Value = 7179.077
$ws.Range("J100").Value = 7021
$ws.Range("L100").Value = 7021
$ws.Range("N100").Value = -8103

$ws.Range("H134").Value = 137750
$ws.Range("J134").Value = 137750
$ws.Range("L134").Value = 137750
$ws.Range("N134").Value = -147890

$ws.Range("H137").Value = 3282.1482
$ws.Range("J137").Value = 3391.3572
$ws.Range("L137").Value = 10174.0716
$ws.Range("N137").Value = -15274.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15010247
$ws.Range("I61").Value = 16677789
$ws.Range("K61").Value = 16677789
$ws.Range("M61").Value = -16677577

$ws.Range("H126").Value = 17864662
$ws.Range("I126").Value = 17864662
$ws.Range("K126").Value = 53593986
$ws.Range("M126").Value = -53591516

$ws.Range("H136").Value = 15010247
$ws.Range("I136").Value = 16677789
$ws.Range("K136").Value = 50033367
$ws.Range("M136").Value = -50030817

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 95974.75
$ws.Range("J57").Value = 95974.75
$ws.Range("L57").Value = 95974.75
$ws.Range("N57").Value = -97414.75

$ws.Range("H87").Value = 94774.5
$ws.Range("J87").Value = 94774.5
$ws.Range("L87").Value = 94774.5
$ws.Range("N87").Value = -97270.5

$ws.Range("H90").Value = 94774.5
$ws.Range("J90").Value = 94774.5
$ws.Range("L90").Value = 284323.5
$ws.Range("N90").Value = -296803.5

$ws.Range("H127").Value = 52950
$ws.Range("J127").Value = 52950
$ws.Range("L127").Value = 52950
$ws.Range("N127").Value = -62870

$ws.Range("H128").Value = 8207.799999999999
$ws.Range("I128").Value = 8207.799999999999
$ws.Range("K128").Value = 24623.4
$ws.Range("M128").Value = -22133.4

$ws.Range("H136").Value = 95974.75
$ws.Range("J136").Value = 95974.75
$ws.Range("L136").Value = 95974.75
$ws.Range("N136").Value = -106174.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3562.1
$ws.Range("J58").Value = 6137.3335
$ws.Range("L58").Value = 6137.3335
$ws.Range("N58").Value = -6543.3335

$ws.Range("H136").Value = 3562.1
$ws.Range("J136").Value = 6137.3335
$ws.Range("L136").Value = 18412.0005
$ws.Range("N136").Value = -23512.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12535.143
$ws.Range("I3").Value = 6893.273
$ws.Range("J3").Value = 33222
$ws.Range("K3").Value = 20679.819
$ws.Range("L3").Value = 99666
$ws.Range("M3").Value = -20567.819
$ws.Range("N3").Value = -99890

$ws.Range("H22").Value = 18416.5
$ws.Range("J22").Value = 33333
$ws.Range("L22").Value = 99999
$ws.Range("N22").Value = -100337

$ws.Range("H25").Value = 14777.667
$ws.Range("I25").Value = 5500
$ws.Range("J25").Value = 33333
$ws.Range("K25").Value = 16500
$ws.Range("L25").Value = 99999
$ws.Range("M25").Value = -16331
$ws.Range("N25").Value = -100337

$ws.Range("H27").Value = 18416.5
$ws.Range("J27").Value = 33333
$ws.Range("L27").Value = 99999
$ws.Range("N27").Value = -100203

$ws.Range("H30").Value = 14777.667
$ws.Range("I30").Value = 5500
$ws.Range("J30").Value = 33333
$ws.Range("K30").Value = 16500
$ws.Range("L30").Value = 99999
$ws.Range("M30").Value = -16398
$ws.Range("N30").Value = -100203

$ws.Range("H112").Value = 17108.25
$ws.Range("I112").Value = 10491.5
$ws.Range("K112").Value = 31474.5
$ws.Range("M112").Value = -30366.5

$ws.Range("H114").Value = 4256.8
$ws.Range("I114").Value = 279.375
$ws.Range("J114").Value = 20166.5
$ws.Range("K114").Value = 838.125
$ws.Range("L114").Value = 60499.5
$ws.Range("M114").Value = 2415.875
$ws.Range("N114").Value = -67007.5

$ws.Range("H127").Value = 500
$ws.Range("J127").Value = 500
$ws.Range("L127").Value = 1500
$ws.Range("N127").Value = -11420

$ws.Range("H138").Value = 14456
$ws.Range("I138").Value = 10685.066
$ws.Range("J138").Value = 33310.668
$ws.Range("K138").Value = 32055.198
$ws.Range("L138").Value = 99932.00399999999
$ws.Range("M138").Value = -26915.198
$ws.Range("N138").Value = -110212.004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14949.167
$ws.Range("I70").Value = 14924.5
$ws.Range("J70").Value = 14998.5
$ws.Range("K70").Value = 14924.5
$ws.Range("L70").Value = 14998.5
$ws.Range("M70").Value = -14654.5
$ws.Range("N70").Value = -15538.5

$ws.Range("H73").Value = 14949.167
$ws.Range("I73").Value = 14924.5
$ws.Range("J73").Value = 14998.5
$ws.Range("K73").Value = 14924.5
$ws.Range("L73").Value = 14998.5
$ws.Range("M73").Value = -13988.5
$ws.Range("N73").Value = -16870.5

$ws.Range("H132").Value = 5787957
$ws.Range("I132").Value = 3058.5173
$ws.Range("J132").Value = 47728468
$ws.Range("K132").Value = 9175.5519
$ws.Range("L132").Value = 143185404
$ws.Range("M132").Value = -6645.5519
$ws.Range("N132").Value = -143190464

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 13716.875
$ws.Range("I107").Value = 13716.875
$ws.Range("K107").Value = 13716.875
$ws.Range("M107").Value = -11796.875

$ws.Range("H137").Value = 116986
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 116986
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 116986
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -127186

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 90655.39999999999
$ws.Range("J56").Value = 110194.25
$ws.Range("L56").Value = 110194.25
$ws.Range("N56").Value = -111622.25
